$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2023-12-08 16:18:55", 0.0014),
    @("2023-12-08 16:19:37", 0.003),
    @("2023-12-08 16:20:01", 0.001)
)

$startRow = 113
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
